$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 447, pushing the existing rows (and everything
# below) down by two. This makes room for a new weekly record.
$ws.Rows.Item(447).Resize(2).Insert()

# Row 447: "Primera" quality record for the new date.
$ws.Cells.Item(447, 1).Value = 8
$ws.Cells.Item(447, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(447, 3).Value = "Coquimbo"
$ws.Cells.Item(447, 4).Value = 45218
$ws.Cells.Item(447, 5).Value = 4
$ws.Cells.Item(447, 6).Value = 100114014
$ws.Cells.Item(447, 7).Value = "Betarraga"
$ws.Cells.Item(447, 8).Value = "Sin especificar"
$ws.Cells.Item(447, 9).Value = "Primera"
$ws.Cells.Item(447, 10).Value = 2000
$ws.Cells.Item(447, 11).Value = 500
$ws.Cells.Item(447, 12).Value = 600
$ws.Cells.Item(447, 13).Value = 550
$ws.Cells.Item(447, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(447, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(447, 16).Value = 183
$ws.Cells.Item(447, 17).Value = 3
$ws.Cells.Item(447, 18).Value = "Hortaliza"

# Row 448: "Segunda" quality record for the new date.
$ws.Cells.Item(448, 1).Value = 8
$ws.Cells.Item(448, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(448, 3).Value = "Coquimbo"
$ws.Cells.Item(448, 4).Value = 45218
$ws.Cells.Item(448, 5).Value = 4
$ws.Cells.Item(448, 6).Value = 100114014
$ws.Cells.Item(448, 7).Value = "Betarraga"
$ws.Cells.Item(448, 8).Value = "Sin especificar"
$ws.Cells.Item(448, 9).Value = "Segunda"
$ws.Cells.Item(448, 10).Value = 1000
$ws.Cells.Item(448, 11).Value = 400
$ws.Cells.Item(448, 12).Value = 450
$ws.Cells.Item(448, 13).Value = 425
$ws.Cells.Item(448, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(448, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(448, 16).Value = 142
$ws.Cells.Item(448, 17).Value = 3
$ws.Cells.Item(448, 18).Value = "Hortaliza"
